{"js": "// This script replaces the date in the title paragraph and the 25\n// multiplication-fact answers in the table with updated values, as a\n// sequence of unique, non-overlapping find/replace operations.\nconst replacements = [\n  [\"2025-06-01 Sunday\", \"2025-06-02 Monday\"],\n  [\"849\u00d76=5094\", \"826\u00d73=2478\"],\n  [\"682\u00d73=2046\", \"259\u00d78=2072\"],\n  [\"294\u00d77=2058\", \"206\u00d78=1648\"],\n  [\"188\u00d78=1504\", \"982\u00d76=5892\"],\n  [\"370\u00d74=1480\", \"319\u00d73=957\"],\n  [\"862\u00d79=7758\", \"837\u00d74=3348\"],\n  [\"559\u00d78=4472\", \"433\u00d76=2598\"],\n  [\"562\u00d78=4496\", \"909\u00d73=2727\"],\n  [\"582\u00d78=4656\", \"626\u00d74=2504\"],\n  [\"847\u00d73=2541\", \"677\u00d73=2031\"],\n  [\"448\u00d78=3584\", \"278\u00d73=834\"],\n  [\"938\u00d73=2814\", \"236\u00d72=472\"],\n  [\"315\u00d76=1890\", \"676\u00d78=5408\"],\n  [\"258\u00d74=1032\", \"542\u00d74=2168\"],\n  [\"766\u00d73=2298\", \"543\u00d72=1086\"],\n  [\"799\u00d79=7191\", \"473\u00d74=1892\"],\n  [\"735\u00d75=3675\", \"778\u00d79=7002\"],\n  [\"713\u00d77=4991\", \"244\u00d79=2196\"],\n  [\"591\u00d79=5319\", \"344\u00d72=688\"],\n  [\"650\u00d73=1950\", \"203\u00d74=812\"],\n  [\"219\u00d79=1971\", \"552\u00d72=1104\"],\n  [\"483\u00d75=2415\", \"724\u00d75=3620\"],\n  [\"218\u00d78=1744\", \"575\u00d74=2300\"],\n  [\"801\u00d78=6408\", \"829\u00d72=1658\"],\n  [\"929\u00d79=8361\", \"507\u00d79=4563\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Updates the date in the title and the 25 multiplication-fact\n# answers in the table to the new values (unique find/replace pairs).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-01 Sunday\", \"2025-06-02 Monday\"),\n    @(\"849\u00d76=5094\", \"826\u00d73=2478\"),\n    @(\"682\u00d73=2046\", \"259\u00d78=2072\"),\n    @(\"294\u00d77=2058\", \"206\u00d78=1648\"),\n    @(\"188\u00d78=1504\", \"982\u00d76=5892\"),\n    @(\"370\u00d74=1480\", \"319\u00d73=957\"),\n    @(\"862\u00d79=7758\", \"837\u00d74=3348\"),\n    @(\"559\u00d78=4472\", \"433\u00d76=2598\"),\n    @(\"562\u00d78=4496\", \"909\u00d73=2727\"),\n    @(\"582\u00d78=4656\", \"626\u00d74=2504\"),\n    @(\"847\u00d73=2541\", \"677\u00d73=2031\"),\n    @(\"448\u00d78=3584\", \"278\u00d73=834\"),\n    @(\"938\u00d73=2814\", \"236\u00d72=472\"),\n    @(\"315\u00d76=1890\", \"676\u00d78=5408\"),\n    @(\"258\u00d74=1032\", \"542\u00d74=2168\"),\n    @(\"766\u00d73=2298\", \"543\u00d72=1086\"),\n    @(\"799\u00d79=7191\", \"473\u00d74=1892\"),\n    @(\"735\u00d75=3675\", \"778\u00d79=7002\"),\n    @(\"713\u00d77=4991\", \"244\u00d79=2196\"),\n    @(\"591\u00d79=5319\", \"344\u00d72=688\"),\n    @(\"650\u00d73=1950\", \"203\u00d74=812\"),\n    @(\"219\u00d79=1971\", \"552\u00d72=1104\"),\n    @(\"483\u00d75=2415\", \"724\u00d75=3620\"),\n    @(\"218\u00d78=1744\", \"575\u00d74=2300\"),\n    @(\"801\u00d78=6408\", \"829\u00d72=1658\"),\n    @(\"929\u00d79=8361\", \"507\u00d79=4563\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
